# Update the build/version string throughout the workbook.
#
# Old version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version string: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet ---------------------------------------------------------
$aboutSheet = $wb.Worksheets.Item("About")

# A2: "Version: <version string>"
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation text containing the version string.
$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Jharia Coal Mine, India, M1686, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet --------------------------------
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S ("build_version") for data rows 2 through 7.
for ($row = 2; $row -le 7; $row++) {
    $dataSheet.Range("S" + $row).Value = $newVersion
}
